# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    3 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 7.143138311642302)
    4 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
    5 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 3.56341032713086)
    6 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
    7 = @(1.455362044514542, 10.34677158129881, 0.7527432677738641, 10.19245300693656, 1, 22.74732990052378)
    8 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    9 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 7.143138311642302)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
